$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.227.87'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").Value = '2.175.96'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.43'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.05%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0910'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.100'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.07%  '

$ws.Range("D14").Value = '2.501.16'
$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.20'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.41%  '

$ws.Range("D16").Value = '2.164.23'
$ws.Range("E16").Value = '  -0.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.764'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.41%  '

$ws.Range("D18").Value = '42.104.61'
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.64'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("B23").Value = 'ImmutableX'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.98%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.18%  '

$ws.Range("E25").Value = '  -0.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.32'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.81%  '

$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.66%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.84%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.67'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +12.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("E33").Value = '  +3.19%  '

$ws.Range("E34").Value = '  -3.53%  '

$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("E36").Value = '  +4.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.54%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0331'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.42%  '

$ws.Range("E41").Value = '  +2.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '58.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.11'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.75%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.31%  '

$ws.Range("B46").Value = 'WOONetwork'
$ws.Range("C46").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.467'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +14.47%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.28%  '

$ws.Range("E48").Value = '  +0.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.28%  '

$ws.Range("E50").Value = '  +0.28%  '

$ws.Range("E51").Value = '  +0.76%  '
